$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newParticipantsQuery = "MATCH (p:participant)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nWITH s, p, samp, f, g, diag`nWHERE g.instrument_model in ['Illumina NovaSeq']`nWITH p`nOPTIONAL MATCH (p)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nWITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID```,`ncoalesce(s.study_name, '') as ``Study Name```,`ncoalesce(s.phs_accession,'') as ``Accession```,`ncoalesce(p.gender,'') as ``Gender```,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY p.participant_id limit 100"

$ws.Cells.Item(2, 2).Value = $newParticipantsQuery

$ws.Range("C3").Select() | Out-Null
